$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update YEAR column (C) from 2020 to 2021 for rows 29-31 and 45-48
$yearRows = @(29, 30, 31, 45, 46, 47, 48)
foreach ($r in $yearRows) {
    $ws.Cells.Item($r, 3).Value = 2021
}

# Update BASE MSRP column (D) values
$ws.Cells.Item(30, 4).Value = 56190
$ws.Cells.Item(45, 4).Value = 76000
$ws.Cells.Item(46, 4).Value = 79250
$ws.Cells.Item(47, 4).Value = 79600
$ws.Cells.Item(48, 4).Value = 82850

# Update the view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
